$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the test name
$ws.Name = "Workflow_1_TestCases"

# Update the active selection
$ws.Range("C13").Select()
